# Add 5 new trial rows (104-108) to Sheet1, continuing the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(103, "Fractals/Version_2/ApoAV-250719-154.jpg", "Reversal", "75_Threat", 0, 0, 0, "right", 3.63,               0.99,               1),
    @(104, "Fractals/Version_2/ApoAV-250719-72.jpg",  "Stable",   "75_Safe",   0, 0, 0, "left",  3.6,                1.43,               1),
    @(105, "Fractals/Version_2/ApoAV-250719-154.jpg", "Reversal", "75_Threat", 0, 0, 1, "right", 4.17,               1.99,               1),
    @(106, "Fractals/Version_2/ApoAV-250719-72.jpg",  "Stable",   "75_Safe",   0, 0, 1, "left",  4.7300000000000004, 1.3699999999999999, 1),
    @(107, "Fractals/Version_2/ApoAV-250719-154.jpg", "Reversal", "75_Threat", 0, 0, 1, "right", 3.64,               2.29,               1)
)

$startRow = 104
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# Update the view state to match the edited workbook (scroll position / active cell).
[void]$ws.Range("A67").Select()
[void]$ws.Range("N104").Select()
